# shipreqdata.xlsx -- "add excel download /mylist, /otherslist"
#
# The shipment-request sheet ("출고요청") gets its single sample data row
# (row 5) replaced by a fresh block of 9 sample rows (rows 5-13), one per
# component, all sharing the same project / developer / department /
# requester columns (A-D) and the same status / empty trailing columns
# (J-K), while E/F/G/H/I vary per part. The "Desc" header cell (B4) is
# restyled to match the yellow header style already used by A4/E4, and the
# active selection moves to I8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 4: give B4 ("Desc") the same yellow header style as A4 ---
$ws.Range("A4").Copy()
$ws.Range("B4").PasteSpecial(-4122)   # xlPasteFormats

# --- Replace the old single sample row with the new 9-row sample block ---
# Delete the existing row 5 first so none of its row-height / per-cell
# style residue (ht="18", s="9" on A5/E5, ...) carries over to the new data.
$ws.Rows(5).Delete()

# H values are written with their original float32 round-trip (the source
# data came through a Java `float`), so the doubles below are the exact
# float32->float64 widenings of 1.213, 0.211, 351.121, 53.2123, 0.1231,
# 0.124, 0.3524, 0.156 and 0.123 respectively.
$rows = @(
    @{ A="TESTPRJA"; B="개발자테스트아이디"; C="전장부품"; D="출고담당자시험용";
       E="2FTS00096A"; F="F6HF2G441AF46";       G="TAIYO YUDEN";            H=1.2130000591278076;  I=530;  J="S1"; K="" },
    @{ A="TESTPRJA"; B="개발자테스트아이디"; C="전장부품"; D="출고담당자시험용";
       E="2FTS00093A"; F="B39242-B4346-P810";    G="EPCOS AG";              H=0.210999995470047;    I=3000; J="S1"; K="" },
    @{ A="TESTPRJA"; B="개발자테스트아이디"; C="전장부품"; D="출고담당자시험용";
       E="2FTC00145A"; F="DLU-2012-25GS1-A1-AT"; G="MAGLAYERS";             H=351.1210021972656;    I=4450; J="S1"; K="" },
    @{ A="TESTPRJA"; B="개발자테스트아이디"; C="전장부품"; D="출고담당자시험용";
       E="2FTP00015A"; F="DPX165950DT-8126A1";   G="TDK";                   H=53.21229934692383;    I=2260; J="S1"; K="" },
    @{ A="TESTPRJA"; B="개발자테스트아이디"; C="전장부품"; D="출고담당자시험용";
       E="2ICZ00186A"; F="RFFM8800TR7";          G="RFMD";                  H=0.12309999763965607;  I=768;  J="S1"; K="" },
    @{ A="TESTPRJA"; B="개발자테스트아이디"; C="전장부품"; D="출고담당자시험용";
       E="2ICT00113A"; F="AR6003XBC2B-R";        G="QUALCOMM INCORPORATED"; H=0.12399999797344208;  I=1275; J="S1"; K="" },
    @{ A="TESTPRJA"; B="개발자테스트아이디"; C="전장부품"; D="출고담당자시험용";
       E="2OSR00012A"; F="1ZCL26000AB0F";        G="KDS";                   H=0.352400004863739;    I=2200; J="S1"; K="" },
    @{ A="TESTPRJA"; B="개발자테스트아이디"; C="전장부품"; D="출고담당자시험용";
       E="2CAC00801A"; F="GRM033R71E471KA01D";   G="MURATA ELEKTRONIK";     H=0.15600000321865082;  I=1560; J="S1"; K="" },
    @{ A="TESTPRJA"; B="개발자테스트아이디"; C="전장부품"; D="출고담당자시험용";
       E="2LL2N5BA11K-R"; F="LQP03TN2N5B02D";    G="MURATA ELEKTRONIK";     H=0.12300000339746475;  I=4000; J="S1"; K="" }
)

$r = 5
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = $row.K
    $r = $r + 1
}

# --- Selection moves to I8 (as recorded by the author's last save) ---
$ws.Range("I8").Select() | Out-Null
